$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A4 from the numeric placeholder 201 to the new string key used as
# the ticket's proto id, adding "TICKET_STAMINA" to the shared strings table.
$ws.Range("A4").Value = "TICKET_STAMINA"

# Widen column A so the new, longer string fits (matches authored width).
$ws.Columns.Item(1).ColumnWidth = 34.86

# Leave the active selection on A4, as in the saved workbook.
$ws.Range("A4").Select()
